$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.988.44"
$ws.Range("E2").Value = "  -1.63%  "

$ws.Range("D3").Value = "3.079.34"
$ws.Range("E3").Value = "  -0.48%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "518.70"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.34%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.17"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.97%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").Value = "3.078.65"
$ws.Range("E8").Value = "  -0.49%  "

$ws.Range("E9").Value = "  +2.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.33"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.26%  "

$ws.Range("E11").Value = "  -1.72%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.399"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.84%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.136"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.82%  "

$ws.Range("D14").Value = "3.613.67"
$ws.Range("E14").Value = "  -0.33%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.23"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.55%  "

$ws.Range("E16").Value = "  -2.31%  "

$ws.Range("D17").Value = "57.126.18"
$ws.Range("E17").Value = "  -1.39%  "

$ws.Range("D18").Value = "3.089.35"
$ws.Range("E18").Value = "  +0.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.88"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.52%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.42"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.82"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.67%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "346.91"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.27%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.78"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.19%  "

$ws.Range("E24").Value = "  -0.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "68.14"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.497"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.82%  "

$ws.Range("E27").Value = "  -2.28%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.997"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.44%  "

$ws.Range("D29").Value = "0.0₃0861"
$ws.Range("E29").Value = "  -6.03%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.17%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.27"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.12%  "

$ws.Range("E32").Value = "  -0.43%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.81"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -9.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.78"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.86%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.88"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +5.52%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "159.50"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.55%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.13"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.00"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "25.54"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.11%  "

$ws.Range("E40").Value = "  -0.82%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0654"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.17%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.04"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.69%  "

$ws.Range("E43").Value = "  +1.32%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.690"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.83%  "

$ws.Range("D45").Value = "2.383.70"
$ws.Range("E45").Value = "  +4.92%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "36.57"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.97%  "

$ws.Range("E47").Value = "  +0.04%  "

$ws.Range("D48").Value = "3.122.70"
$ws.Range("E48").Value = "  -0.36%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0263"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.23%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.953"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.52%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.93"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.67%  "
